# Changes to fighting matrices
# Updates the "Blaster" and "Bubbler" rows of the fighting profile table.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Blaster row
Replace-Text "SR+1" "RS+1"
Replace-Text "SR+5" "PS+4"
Replace-Text "RS+3" "RW+3"
Replace-Text "CP+0" "PW+1"

# Bubbler row
Replace-Text "MP+2" "PW+2"
Replace-Text "PR+3" "PR+2"
Replace-Text "RS+2" "RW+2"
Replace-Text "PW+0" "CP+1"
